# Slide 17: fix typo in the "delete a branch" bullet.
# "git branch (-d) <name>" -> "git branch -d <name>"
#
# The line is built from several runs; the two runs that hold the
# "(-" and "d) <name>" text need to become "-d " and "<name>"
# respectively, while every other run (and each run's formatting)
# stays untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(17)
$shape = $s.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange

$full = $tr.Text
$anchor = $full.IndexOf("To delete a branch")
$openParenIdx = $full.IndexOf("(-", $anchor)

if ($openParenIdx -lt 0) {
    throw "Could not locate '(-' run after 'To delete a branch' on slide 17"
}

# "(-" is 2 characters, followed immediately by "d) <name>"
$runOpenParen = $tr.Characters($openParenIdx + 1, 2)
$runDName = $tr.Characters($openParenIdx + 3, "d) <name>".Length)

# Update the second run first so the first run's offsets stay valid.
$runDName.Text = "<name>"
$runOpenParen.Text = "-d "
